# Refresh the cryptocurrency Price (D) / Volume(1h) (E) figures with the latest pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of "Price" cells (column D) hold purely numeric-looking text (e.g. "377.66").
# Assigning such a string straight to .Value would make Excel auto-convert it to a real
# number, silently dropping the original text formatting (trailing zeros, etc.). Force
# just those specific cells to Text format first so the refreshed values stay literal
# strings, matching how this column is already stored in the sheet.
$forceTextRefs = @("D5", "D6", "D8", "D12", "D13", "D17", "D19", "D20", "D21", "D23", "D24", "D27", "D29", "D30", "D31", "D33", "D35", "D36", "D40", "D44", "D45")
foreach ($ref in $forceTextRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range('D2').Value = '51.391.95'
$ws.Range('E2').Value = '  -1.35%  '
$ws.Range('D3').Value = '2.924.57'
$ws.Range('E3').Value = '  -2.52%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '377.66'
$ws.Range('E5').Value = '  +6.85%  '
$ws.Range('D6').Value = '103.05'
$ws.Range('E6').Value = '  -3.31%  '
$ws.Range('E7').Value = '  -2.45%  '
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  -0.24%  '
$ws.Range('E9').Value = '  -3.50%  '
$ws.Range('E10').Value = '  -2.61%  '
$ws.Range('E11').Value = '  -0.22%  '
$ws.Range('D12').Value = '0.0835'
$ws.Range('E12').Value = '  -2.21%  '
$ws.Range('D13').Value = '18.32'
$ws.Range('E13').Value = '  -3.58%  '
$ws.Range('D14').Value = '3.385.32'
$ws.Range('E14').Value = '  -2.51%  '
$ws.Range('E15').Value = '  -3.37%  '
$ws.Range('D16').Value = '2.919.66'
$ws.Range('E16').Value = '  -2.59%  '
$ws.Range('D17').Value = '0.929'
$ws.Range('E17').Value = '  -8.37%  '
$ws.Range('D18').Value = '51.330.73'
$ws.Range('E18').Value = '  -1.52%  '
$ws.Range('D19').Value = '3.42'
$ws.Range('E19').Value = '  -0.08%  '
$ws.Range('D20').Value = '7.37'
$ws.Range('E20').Value = '  -1.35%  '
$ws.Range('D21').Value = '12.95'
$ws.Range('E21').Value = '  -4.11%  '
$ws.Range('E22').Value = '  -2.35%  '
$ws.Range('D23').Value = '68.36'
$ws.Range('E23').Value = '  -0.94%  '
$ws.Range('D24').Value = '261.26'
$ws.Range('E24').Value = '  -1.02%  '
$ws.Range('E25').Value = '  +1.88%  '
$ws.Range('E26').Value = '  -3.54%  '
$ws.Range('D27').Value = '4.13'
$ws.Range('E27').Value = '  -4.75%  '
$ws.Range('E28').Value = '  +0.04%  '
$ws.Range('D29').Value = '25.67'
$ws.Range('E29').Value = '  -4.26%  '
$ws.Range('D30').Value = '7.21'
$ws.Range('E30').Value = '  -2.94%  '
$ws.Range('D31').Value = '6.83'
$ws.Range('E31').Value = '  +7.08%  '
$ws.Range('E32').Value = '  -4.59%  '
$ws.Range('D33').Value = '9.81'
$ws.Range('E33').Value = '  -3.91%  '
$ws.Range('E34').Value = '  -3.60%  '
$ws.Range('D35').Value = '51.05'
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('D36').Value = '34.08'
$ws.Range('E36').Value = '  -4.77%  '
$ws.Range('E37').Value = '  +0.30%  '
$ws.Range('E38').Value = '  -3.85%  '
$ws.Range('E39').Value = '  -8.68%  '
$ws.Range('D40').Value = '16.93'
$ws.Range('E40').Value = '  -3.68%  '
$ws.Range('E41').Value = '  -9.79%  '
$ws.Range('E42').Value = '  -7.79%  '
$ws.Range('E43').Value = '  -2.33%  '
$ws.Range('D44').Value = '123.14'
$ws.Range('E44').Value = '  -1.04%  '
$ws.Range('D45').Value = '21.56'
$ws.Range('E45').Value = '  -5.52%  '
$ws.Range('E46').Value = '  -2.89%  '
$ws.Range('E47').Value = '  +11.97%  '
$ws.Range('D48').Value = '2.027.23'
$ws.Range('E48').Value = '  -4.29%  '
$ws.Range('E49').Value = '  -2.35%  '
$ws.Range('E50').Value = '  -4.73%  '
$ws.Range('D51').Value = '3.202.39'
$ws.Range('E51').Value = '  -2.83%  '
